$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column widths: col A -> 14.7109375, col B -> 15.42578125
# (ColumnWidth is quantized internally to pixel boundaries by the COM width
# model, so we pick the character-width inputs that land closest to the
# exact target stored widths.)
$ws.Columns.Item(1).ColumnWidth = 13.833333333333334
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

# Update values in existing rows 1-4
$ws.Range("A1").Value = 0.024381441209963044
$ws.Range("B1").Value = -0.024381441212844003

$ws.Range("A2").Value = -0.014453546071601079
$ws.Range("B2").Value = 0.014453546065520741

$ws.Range("A3").Value = -0.014857821428728784
$ws.Range("B3").Value = 0.014857821421021473

$ws.Range("A4").Value = 0.0035965660576935698
$ws.Range("B4").Value = -0.0035965660608587788

# Add new row 5
$ws.Range("A5").Value = 0.064024928246492807
$ws.Range("B5").Value = -0.064024928251123908
